# Auto update stock data
# Update the "latest" row (2025/12/18 -> 2025/12/19) date and EBITDA figures
# for each company block in the sheet. All cells in columns A and B are
# stored as text in the workbook, so we force text formatting before
# writing the new values to avoid Excel auto-converting them to a date
# serial number / floating point number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Value
    )
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
}

# Row 2 - Alcoa
Set-TextValue "A2" "2025/12/19"
Set-TextValue "B2" "5.92"

# Row 8 - Rio Tinto
Set-TextValue "A8" "2025/12/19"
Set-TextValue "B8" "8.24"

# Row 14 - Norsk Hydro
Set-TextValue "A14" "2025/12/19"
Set-TextValue "B14" "2.99"

# Row 20 - Reliance
Set-TextValue "A20" "2025/12/19"
Set-TextValue "B20" "13.00"

# Row 26 - Kaiser
Set-TextValue "A26" "2025/12/19"
Set-TextValue "B26" "10.98"

# Row 32 - Ryerson
Set-TextValue "A32" "2025/12/19"
Set-TextValue "B32" "27.66"

# Row 38 - Alro Steel (date only, EBITDA unchanged)
Set-TextValue "A38" "2025/12/19"

# Row 44 - Ultra
Set-TextValue "A44" "2025/12/19"
Set-TextValue "B44" "10.94"

# Row 50 - Benchmark
Set-TextValue "A50" "2025/12/19"
Set-TextValue "B50" "11.36"

# Row 56 - Celestica
Set-TextValue "A56" "2025/12/19"
Set-TextValue "B56" "28.44"

# Row 62 - Jabil
Set-TextValue "A62" "2025/12/19"
Set-TextValue "B62" "10.88"

# Row 68 - Flex
Set-TextValue "A68" "2025/12/19"
Set-TextValue "B68" "13.12"

# Row 74 - MKS
Set-TextValue "A74" "2025/12/19"
Set-TextValue "B74" "16.22"
